$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long numeric-looking account number that must stay TEXT
# (matches the original cell's storage type). Force text entry, then
# restore the cell's original number format/style (copied from a
# neighboring cell that already carries that exact style) so the
# cell style index is unchanged.
$b3 = $ws.Range("B3")
$b3.NumberFormat = "@"
$b3.Value = "2570314725427075"
$ws.Range("C3").Copy() | Out-Null
$b3.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 09.05.2024"

# Row 6
$ws.Range("B6").Value = "10.05."
$ws.Range("C6").Value = "11.05."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 6328499"
$ws.Range("E6").Value = "40,40-"

# Row 7
$ws.Range("B7").Value = "13.05."
$ws.Range("C7").Value = "14.05."
$ws.Range("D7").Value = "ZALANDO MKTPLC EU WTBGCS"
$ws.Range("E7").Value = "130,00-"

# Row 8
$ws.Range("B8").Value = "14.05."
$ws.Range("C8").Value = "15.05."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,85-"

# Row 9 - the transaction that used to live here is gone; the row is
# blanked out. E9 switches from the right-aligned amount style to the
# centered blank style (index 13 in the original style table) so set
# the alignment explicitly rather than copying a format.
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$e9 = $ws.Range("E9")
$e9.Value = ""
$e9.HorizontalAlignment = -4108  # xlCenter
$e9.VerticalAlignment = -4108    # xlCenter
$e9.WrapText = $true

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 18.05.2024"
$ws.Range("E12").Value = "195,25-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 28.05.2024"
